$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ClinicianLanding")

# Row 19 text updates
$ws.Range("B19").Value = "Compare the patient name"
$ws.Range("C19").Value = "Compare"

# F19: new patient-name value with its own bold/Open Sans styling
$f19 = $ws.Range("F19")
$f19.Value = "SM2ITHEE, Sophia"
$f19.Borders.LineStyle = -4142
$font = $f19.Font
$font.Name = "Open Sans"
$font.Bold = $true
$font.Size = 9
$font.Color = 3552822

# Row 19 grew slightly taller once the new text/font were added
$ws.Rows.Item(19).RowHeight = 15.75

# Selection moved to F19 after the edit
[void]$ws.Range("F19").Select()
